$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45203
$ws.Range("J2").Value = 100
# Row 3
$ws.Range("D3").Value = 45203
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1500
$ws.Range("P3").Value = 500
# Row 4
$ws.Range("D4").Value = 45219
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2250
$ws.Range("P4").Value = 750
# Row 5
$ws.Range("D5").Value = 44838
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1250
$ws.Range("P5").Value = 417
# Row 6
$ws.Range("D6").Value = 44838
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 1000
$ws.Range("P6").Value = 333
# Row 7
$ws.Range("D7").Value = 45163
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2500
$ws.Range("P7").Value = 833
# Row 8
$ws.Range("D8").Value = 45215
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 2000
$ws.Range("P8").Value = 667
# Row 9
$ws.Range("D9").Value = 45135
$ws.Range("J9").Value = 70
# Row 10
$ws.Range("D10").Value = 45175
$ws.Range("J10").Value = 150
# Row 11
$ws.Range("D11").Value = 45134
$ws.Range("J11").Value = 50
# Row 12
$ws.Range("D12").Value = 45145
$ws.Range("J12").Value = 60
# Row 13
$ws.Range("D13").Value = 45145
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 80
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 2000
$ws.Range("P13").Value = 667
# Row 14
$ws.Range("D14").Value = 44832
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 1300
$ws.Range("M14").Value = 1250
$ws.Range("P14").Value = 417
# Row 15
$ws.Range("D15").Value = 44832
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 150
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 1000
$ws.Range("P15").Value = 333
# Row 16
$ws.Range("D16").Value = 45176
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("P16").Value = 833
# Row 17
$ws.Range("D17").Value = 45195
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2500
$ws.Range("P17").Value = 833
# Row 19
$ws.Range("D19").Value = 45161
# Row 20
$ws.Range("D20").Value = 45160
$ws.Range("J20").Value = 100
# Row 21
$ws.Range("D21").Value = 45173
# Row 23
$ws.Range("D23").Value = 45146
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 2500
$ws.Range("L23").Value = 2500
$ws.Range("M23").Value = 2500
$ws.Range("P23").Value = 833
# Row 24
$ws.Range("D24").Value = 45146
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = 2000
$ws.Range("P24").Value = 667
# Row 25
$ws.Range("D25").Value = 44846
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 1200
$ws.Range("L25").Value = 1300
$ws.Range("M25").Value = 1250
$ws.Range("P25").Value = 417
# Row 26
$ws.Range("D26").Value = 44846
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 1000
$ws.Range("P26").Value = 333
# Row 27
$ws.Range("D27").Value = 45191
$ws.Range("J27").Value = 100
# Row 30
$ws.Range("D30").Value = 45166
$ws.Range("J30").Value = 120
# Row 31
$ws.Range("D31").Value = 45133
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = 2500
$ws.Range("P31").Value = 833
# Row 32
$ws.Range("D32").Value = 45149
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 2500
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = 2500
$ws.Range("P32").Value = 833
# Row 33
$ws.Range("D33").Value = 45149
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 80
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = 2000
$ws.Range("P33").Value = 667
